# Remove the "TestResultExcelFilePath" output-file-path column (column H) from
# the NI Scenario input sheets. The two sheets ("ProcessPayrollForNIWeekly"
# and "TestReports") were edited together (grouped), with
# "ProcessPayrollForNIWeekly" ending up as the active sheet/tab, and column H
# fully selected (whole-column selection) on both sheets after the delete.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProcessPayrollForNIWeekly")
$wsReport = $wb.Worksheets.Item("TestReports")

# Delete the whole "TestResultExcelFilePath" column (H) on both sheets -
# this removes the header cell and every data-row value in that column and
# shifts the columns to its right one place to the left.
$wsInput.Columns("H:H").Delete()
$wsReport.Columns("H:H").Delete()

# "ProcessPayrollForNIWeekly" becomes the active sheet, with the (now empty)
# former column H left selected as a whole column.
[void]$wsInput.Activate()
[void]$wsInput.Range("H1:H1048576").Select()

[void]$wsReport.Range("H1:H1048576").Select()
[void]$wsInput.Activate()
